# Commit: "I0 and IF added"
# Adds two new columns to Sheet1:
#   I1 = "I0"  (header, same style as existing headers)
#   J1 = "IF"  (header, same style as existing headers)
#   I2:I51 / J2:J51 = per-row numeric data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1/J1 -- copy the formatting of an existing header cell (H1)
# so the new headers pick up the same bold/centered/bordered style already
# used by the rest of row 1, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data values for I2:J51 (column I = "I0", column J = "IF")
$values = @(
    @(8,8),
    @(6,6),
    @(7,7),
    @(5,6),
    @(6,6),
    @(6,7),
    @(7,7),
    @(7,7),
    @(5,5),
    @(6,6),
    @(7,7),
    @(5,5),
    @(6,6),
    @(6,6),
    @(7,7),
    @(6,7),
    @(8,8),
    @(6,6),
    @(8,8),
    @(6,6),
    @(5,6),
    @(7,7),
    @(6,6),
    @(6,6),
    @(8,8),
    @(7,7),
    @(7,7),
    @(8,8),
    @(8,8),
    @(7,7),
    @(6,6),
    @(7,7),
    @(6,6),
    @(3,4),
    @(7,7),
    @(8,8),
    @(7,7),
    @(7,7),
    @(7,7),
    @(7,7),
    @(6,7),
    @(7,7),
    @(8,8),
    @(8,8),
    @(9,9),
    @(7,8),
    @(5,6),
    @(6,6),
    @(7,7),
    @(8,8)
)

for ($r = 0; $r -lt $values.Length; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value  = $values[$r][0]   # column I
    $ws.Cells.Item($row, 10).Value = $values[$r][1]   # column J
}
